# Auto-generated edit script applying Kraken_Profits market-data refresh
# Updates currentAveragePrice* and computed Leve profit columns (H, I, J, K, L, M, N)
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 290.9375
$ws.Range("I15").Value = 290.9375
$ws.Range("K15").Value = 872.8125
$ws.Range("M15").Value = -703.8125
$ws.Range("H32").Value = 10300
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 10300
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 10300
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -10952
$ws.Range("H33").Value = 192
$ws.Range("I33").Value = 192
$ws.Range("K33").Value = 192
$ws.Range("M33").Value = 37
$ws.Range("H51").Value = 6600
$ws.Range("H64").Value = 8001.5
$ws.Range("J64").Value = 8001.5
$ws.Range("L64").Value = 8001.5
$ws.Range("N64").Value = -8497.5
$ws.Range("H67").Value = 8001.5
$ws.Range("J67").Value = 8001.5
$ws.Range("L67").Value = 8001.5
$ws.Range("N67").Value = -9717.5
$ws.Range("H69").Value = 34666.668
$ws.Range("J69").Value = 2000
$ws.Range("L69").Value = 6000
$ws.Range("N69").Value = -7748
$ws.Range("H72").Value = 34666.668
$ws.Range("J72").Value = 2000
$ws.Range("L72").Value = 18000
$ws.Range("N72").Value = -26736
$ws.Range("H111").Value = 294.33334
$ws.Range("I111").Value = 294.33334
$ws.Range("K111").Value = 883.0000200000001
$ws.Range("M111").Value = 2183.99998
$ws.Range("H132").Value = 5439.4287
$ws.Range("I132").Value = 4679.3335
$ws.Range("K132").Value = 14038.0005
$ws.Range("M132").Value = -11508.0005
$ws.Range("H135").Value = 1669.5
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2540
$ws.Range("I97").Value = 2098.6428
$ws.Range("J97").Value = 4599.6665
$ws.Range("K97").Value = 2098.6428
$ws.Range("L97").Value = 4599.6665
$ws.Range("M97").Value = -1602.6428
$ws.Range("N97").Value = -5591.6665
$ws.Range("H110").Value = 998.5
$ws.Range("I110").Value = 998.5
$ws.Range("K110").Value = 998.5
$ws.Range("M110").Value = 1046.5
$ws.Range("H132").Value = 2078.4614
$ws.Range("I132").Value = 1370.1666
$ws.Range("J132").Value = 2685.5715
$ws.Range("K132").Value = 4110.4998
$ws.Range("L132").Value = 8056.7145
$ws.Range("M132").Value = -1580.4998
$ws.Range("N132").Value = -13116.7145

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 99995
$ws.Range("J59").Value = 99995
$ws.Range("L59").Value = 99995
$ws.Range("N59").Value = -101689
$ws.Range("H80").Value = 551.75
$ws.Range("I80").Value = 350
$ws.Range("J80").Value = 753.5
$ws.Range("K80").Value = 350
$ws.Range("L80").Value = 753.5
$ws.Range("M80").Value = 648
$ws.Range("N80").Value = -2749.5
$ws.Range("H83").Value = 551.75
$ws.Range("I83").Value = 350
$ws.Range("J83").Value = 753.5
$ws.Range("K83").Value = 1750
$ws.Range("L83").Value = 3767.5
$ws.Range("M83").Value = 3242
$ws.Range("N83").Value = -13751.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3218.4
$ws.Range("I22").Value = 3773
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 3773
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -3423
$ws.Range("N22").Value = -1700
$ws.Range("H31").Value = 7523.4707
$ws.Range("I31").Value = 7108.1665
$ws.Range("K31").Value = 7108.1665
$ws.Range("M31").Value = -6813.1665
$ws.Range("H34").Value = 7523.4707
$ws.Range("I34").Value = 7108.1665
$ws.Range("K34").Value = 7108.1665
$ws.Range("M34").Value = -6906.1665
$ws.Range("H57").Value = 6000
$ws.Range("J57").Value = 6000
$ws.Range("L57").Value = 6000
$ws.Range("N57").Value = -7120
$ws.Range("H132").Value = 2907
$ws.Range("I132").Value = 2900
$ws.Range("K132").Value = 8700
$ws.Range("M132").Value = -6170
$ws.Range("H134").Value = 5981.2
$ws.Range("I134").Value = 5981.2
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 17943.6
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -15408.6
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 99995
$ws.Range("J137").Value = 99995
$ws.Range("L137").Value = 99995
$ws.Range("N137").Value = -110195

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1902.4
$ws.Range("I12").Value = 307.33334
$ws.Range("J12").Value = 2586
$ws.Range("K12").Value = 922.0000200000001
$ws.Range("L12").Value = 7758
$ws.Range("M12").Value = -749.0000200000001
$ws.Range("N12").Value = -8104
$ws.Range("H48").Value = 1000
$ws.Range("J48").Value = 1000
$ws.Range("L48").Value = 3000
$ws.Range("N48").Value = -3500
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H55").Value = 1036.875
$ws.Range("I55").Value = 1036.875
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 3110.625
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -2933.625
$ws.Range("N55").ClearContents()
$ws.Range("H113").Value = 723
$ws.Range("I113").Value = 888
$ws.Range("J113").Value = 681.75
$ws.Range("K113").Value = 2664
$ws.Range("L113").Value = 2045.25
$ws.Range("M113").Value = -494
$ws.Range("N113").Value = -6385.25

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5284
$ws.Range("I126").Value = 5284
$ws.Range("K126").Value = 15852
$ws.Range("M126").Value = -13382

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1266.6666
$ws.Range("J2").Value = 1266.6666
$ws.Range("L2").Value = 1266.6666
$ws.Range("N2").Value = -1490.6666
$ws.Range("H7").Value = 2897.4
$ws.Range("I7").Value = 2897.4
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2897.4
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2785.4
$ws.Range("N7").ClearContents()
$ws.Range("H32").Value = 2299
$ws.Range("I32").Value = 2299
$ws.Range("K32").Value = 2299
$ws.Range("M32").Value = -1982
$ws.Range("H38").Value = 30000
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H40").Value = 5667
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H81").Value = 40181
$ws.Range("J81").Value = 40181
$ws.Range("L81").Value = 40181
$ws.Range("N81").Value = -42177
$ws.Range("H82").Value = 1420.2307
$ws.Range("J82").Value = 1441.25
$ws.Range("L82").Value = 1441.25
$ws.Range("N82").Value = -2163.25
$ws.Range("H84").Value = 40181
$ws.Range("J84").Value = 40181
$ws.Range("L84").Value = 120543
$ws.Range("N84").Value = -130527
$ws.Range("H85").Value = 1420.2307
$ws.Range("J85").Value = 1441.25
$ws.Range("L85").Value = 1441.25
$ws.Range("N85").Value = -3937.25
$ws.Range("H126").Value = 2897.4
$ws.Range("I126").Value = 2897.4
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8692.200000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6222.200000000001
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 4158.6
$ws.Range("I132").Value = 4824.75
$ws.Range("K132").Value = 14474.25
$ws.Range("M132").Value = -11944.25

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 13116.167
$ws.Range("I3").Value = 30001
$ws.Range("J3").Value = 9739.200000000001
$ws.Range("K3").Value = 30001
$ws.Range("L3").Value = 9739.200000000001
$ws.Range("M3").Value = -29887
$ws.Range("N3").Value = -9967.200000000001
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H11").Value = 7833
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H45").Value = 39069
$ws.Range("J45").Value = 44888
$ws.Range("L45").Value = 44888
$ws.Range("N45").Value = -45870
$ws.Range("H52").Value = 10015021
$ws.Range("I52").Value = 10015021
$ws.Range("K52").Value = 10015021
$ws.Range("M52").Value = -10014795
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H122").Value = 5417.1665
$ws.Range("I122").Value = 4899.6
$ws.Range("J122").Value = 8005
$ws.Range("K122").Value = 14698.8
$ws.Range("L122").Value = 24015
$ws.Range("M122").Value = -12248.8
$ws.Range("N122").Value = -28915
$ws.Range("H132").Value = 11038.5
$ws.Range("I132").Value = 9332.75
$ws.Range("J132").Value = 14450
$ws.Range("K132").Value = 27998.25
$ws.Range("L132").Value = 43350
$ws.Range("M132").Value = -25468.25
$ws.Range("N132").Value = -48410
$ws.Range("H136").Value = 1220.1538
$ws.Range("I136").Value = 1220.1538
$ws.Range("K136").Value = 3660.4614
$ws.Range("M136").Value = -1110.4614
